{"js": "// The document had two <id>...</id> tags whose opening tag, inner value,\n// and closing tag were each split across separate runs (with the\n// inner-value run carrying lighter/plain formatting than the tag runs).\n// The edit collapses each of those three runs into a single run (using\n// the opening tag's run formatting), and also renames the second id's\n// value from \"p127r_a1\" to \"p127r_1\".\n//\n// Using search() + insertText(..., \"Replace\") on the whole \"<id>...</id>\"\n// span naturally merges the matched runs into one run (seeded from the\n// first run's formatting), which reproduces the diff.\n\nconst body = context.document.body;\n\nconst firstMatches = body.search(\"<id>p126v_3</id>\", { matchCase: true });\nfirstMatches.load(\"items\");\nconst secondMatches = body.search(\"<id>p127r_a1</id>\", { matchCase: true });\nsecondMatches.load(\"items\");\n\nawait context.sync();\n\nif (firstMatches.items.length > 0) {\n  firstMatches.items[0].insertText(\"<id>p126v_3</id>\", \"Replace\");\n}\nif (secondMatches.items.length > 0) {\n  secondMatches.items[0].insertText(\"<id>p127r_1</id>\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document had two <id>...</id> tags whose opening tag, inner value,\n# and closing tag each lived in their own run (the inner-value run using\n# lighter/plain formatting vs. the tag runs' Courier New / 7f6000 look).\n# The edit collapses each of those three runs into a single run carrying\n# the opening tag's run formatting, and also renames the second id's\n# value from \"p127r_a1\" to \"p127r_1\".\n#\n# Find.Execute with ReplaceWith/Replace (wdReplaceOne) rewrites the whole\n# matched \"<id>...</id>\" span as one run seeded from the first matched\n# run's formatting - exactly reproducing the merge.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 2\n\n$range1 = $d.Content\n$range1.Find.Execute(\n    \"<id>p126v_3</id>\",\n    $false,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    $wdFindContinue,\n    $false,\n    \"<id>p126v_3</id>\",\n    $wdReplaceOne\n)\n\n$range2 = $d.Content\n$range2.Find.Execute(\n    \"<id>p127r_a1</id>\",\n    $false,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    $wdFindContinue,\n    $false,\n    \"<id>p127r_1</id>\",\n    $wdReplaceOne\n)\n"}
